# Auto-generated script to update market-price derived columns (H:N)
# across several leve-profit worksheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4455.909
$ws.Range("I28").Value = 252.66667
$ws.Range("J28").Value = 9499.799999999999
$ws.Range("K28").Value = 252.66667
$ws.Range("L28").Value = 9499.799999999999
$ws.Range("M28").Value = 232.33333
$ws.Range("N28").Value = -10469.8
$ws.Range("H33").Value = 142.76923
$ws.Range("I33").Value = 195.25
$ws.Range("J33").Value = 119.44444
$ws.Range("K33").Value = 195.25
$ws.Range("L33").Value = 119.44444
$ws.Range("M33").Value = 33.75
$ws.Range("N33").Value = -577.44444
$ws.Range("H125").Value = 641.7
$ws.Range("I125").Value = 766.4
$ws.Range("J125").Value = 517
$ws.Range("K125").Value = 6897.599999999999
$ws.Range("L125").Value = 4653
$ws.Range("M125").Value = -4437.599999999999
$ws.Range("N125").Value = -9573
$ws.Range("H126").Value = 40000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 40000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880
$ws.Range("H127").Value = 2855
$ws.Range("I127").Value = 2855
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 8565
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -3605
$ws.Range("H128").Value = 40000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 40000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960
$ws.Range("H129").Value = 925.4789
$ws.Range("I129").Value = 935
$ws.Range("J129").Value = 924.6
$ws.Range("K129").Value = 2805
$ws.Range("L129").Value = 2773.8
$ws.Range("M129").Value = 2195
$ws.Range("N129").Value = -12773.8
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 2680.6875
$ws.Range("I131").Value = 444.9091
$ws.Range("J131").Value = 7599.4
$ws.Range("K131").Value = 1334.7273
$ws.Range("L131").Value = 22798.2
$ws.Range("M131").Value = 3705.2727
$ws.Range("N131").Value = -32878.2
$ws.Range("H132").Value = 1188.9642
$ws.Range("I132").Value = 1196.2307
$ws.Range("J132").Value = 1094.5
$ws.Range("K132").Value = 3588.6921
$ws.Range("L132").Value = 3283.5
$ws.Range("M132").Value = -1058.6921
$ws.Range("N132").Value = -8343.5
$ws.Range("H133").Value = 89000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 89000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120
$ws.Range("H134").Value = 42166.668
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 42166.668
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 42166.668
$ws.Range("N134").Value = -52306.668
$ws.Range("H135").Value = 438.5
$ws.Range("I135").Value = 387.86667
$ws.Range("J135").Value = 691.6667
$ws.Range("K135").Value = 3490.80003
$ws.Range("L135").Value = 6225.0003
$ws.Range("M135").Value = -955.8000299999999
$ws.Range("N135").Value = -11295.0003
$ws.Range("H136").Value = 65000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 65000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 65000
$ws.Range("N136").Value = -75200
$ws.Range("H137").Value = 2053.238
$ws.Range("I137").Value = 1995
$ws.Range("J137").Value = 2056.15
$ws.Range("K137").Value = 5985
$ws.Range("L137").Value = 6168.450000000001
$ws.Range("M137").Value = -3435
$ws.Range("N137").Value = -11268.45
$ws.Range("H138").Value = 4663.1313
$ws.Range("I138").Value = 4485.2666
$ws.Range("J138").Value = 4721.1304
$ws.Range("K138").Value = 13455.7998
$ws.Range("L138").Value = 14163.3912
$ws.Range("M138").Value = -8315.799800000001
$ws.Range("N138").Value = -24443.3912
$ws.Range("H139").Value = 48100
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 48100
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 48100
$ws.Range("N139").Value = -58380
$ws.Range("H140").Value = 163932.83
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 163932.83
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 163932.83
$ws.Range("N140").Value = -174292.83
$ws.Range("H141").Value = 2003527.2
$ws.Range("I141").Value = 5601735
$ws.Range("J141").Value = 4522.778
$ws.Range("K141").Value = 16805205
$ws.Range("L141").Value = 13568.334
$ws.Range("M141").Value = -16800025
$ws.Range("N141").Value = -23928.334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1110488.5
$ws.Range("I2").Value = 1165763
$ws.Range("J2").Value = 4999
$ws.Range("K2").Value = 1165763
$ws.Range("L2").Value = 4999
$ws.Range("M2").Value = -1165650
$ws.Range("N2").Value = -5225
$ws.Range("H32").Value = 2596.9678
$ws.Range("I32").Value = 2035.7976
$ws.Range("K32").Value = 2035.7976
$ws.Range("M32").Value = -1748.7976
$ws.Range("H110").Value = 1051.125
$ws.Range("I110").Value = 1051.125
$ws.Range("K110").Value = 1051.125
$ws.Range("M110").Value = 993.875
$ws.Range("H116").Value = 1110488.5
$ws.Range("I116").Value = 1165763
$ws.Range("J116").Value = 4999
$ws.Range("K116").Value = 1165763
$ws.Range("L116").Value = 4999
$ws.Range("M116").Value = -1163469
$ws.Range("N116").Value = -9587
$ws.Range("H132").Value = 2413.6428
$ws.Range("I132").Value = 1834.24
$ws.Range("J132").Value = 3265.7058
$ws.Range("K132").Value = 5502.72
$ws.Range("L132").Value = 9797.117400000001
$ws.Range("M132").Value = -2972.72
$ws.Range("N132").Value = -14857.1174

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1110488.5
$ws.Range("I3").Value = 1165763
$ws.Range("J3").Value = 4999
$ws.Range("K3").Value = 1165763
$ws.Range("L3").Value = 4999
$ws.Range("M3").Value = -1165649
$ws.Range("N3").Value = -5227
$ws.Range("H86").Value = 2143.5264
$ws.Range("I86").Value = 1948.7142
$ws.Range("J86").Value = 2689
$ws.Range("K86").Value = 1948.7142
$ws.Range("L86").Value = 2689
$ws.Range("M86").Value = -825.7141999999999
$ws.Range("N86").Value = -4935
$ws.Range("H89").Value = 2143.5264
$ws.Range("I89").Value = 1948.7142
$ws.Range("J89").Value = 2689
$ws.Range("K89").Value = 9743.571
$ws.Range("L89").Value = 13445
$ws.Range("M89").Value = -4127.571
$ws.Range("N89").Value = -24677
$ws.Range("H99").Value = 1181.6
$ws.Range("I99").Value = 1181.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1181.6
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 316.4000000000001
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 2319.55
$ws.Range("I107").Value = 1716.2307
$ws.Range("K107").Value = 1716.2307
$ws.Range("M107").Value = 203.7692999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1593.5
$ws.Range("I99").Value = 1381
$ws.Range("J99").Value = 1735.1666
$ws.Range("K99").Value = 1381
$ws.Range("L99").Value = 1735.1666
$ws.Range("M99").Value = 117
$ws.Range("N99").Value = -4731.1666
$ws.Range("H122").Value = 1245.75
$ws.Range("I122").Value = 1245.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3737.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1287.25
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1593.5
$ws.Range("I126").Value = 1381
$ws.Range("J126").Value = 1735.1666
$ws.Range("K126").Value = 4143
$ws.Range("L126").Value = 5205.4998
$ws.Range("M126").Value = -1673
$ws.Range("N126").Value = -10145.4998
$ws.Range("H134").Value = 2724.0386
$ws.Range("I134").Value = 1615.4762
$ws.Range("J134").Value = 7380
$ws.Range("K134").Value = 4846.4286
$ws.Range("L134").Value = 22140
$ws.Range("M134").Value = -2311.4286
$ws.Range("N134").Value = -27210

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 111189720
$ws.Range("I50").Value = 175018.5
$ws.Range("J50").Value = 200001490
$ws.Range("K50").Value = 525055.5
$ws.Range("L50").Value = 600004470
$ws.Range("M50").Value = -524574.5
$ws.Range("N50").Value = -600005432
$ws.Range("H53").Value = 111189720
$ws.Range("I53").Value = 175018.5
$ws.Range("J53").Value = 200001490
$ws.Range("K53").Value = 525055.5
$ws.Range("L53").Value = 600004470
$ws.Range("M53").Value = -524574.5
$ws.Range("N53").Value = -600005432
$ws.Range("H129").Value = 30047.24
$ws.Range("J129").Value = 41503.223
$ws.Range("L129").Value = 124509.669
$ws.Range("N129").Value = -134509.669

